# Auto-generated edit script: refresh market price data (columns H-N) across multiple sheets
# This mirrors a scheduled-runner update of Kujata_Profits market price columns.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3205631.8
$ws.Range("I18").Value = 255
$ws.Range("J18").Value = 13890221
$ws.Range("K18").Value = 255
$ws.Range("L18").Value = 13890221
$ws.Range("M18").Value = 29
$ws.Range("N18").Value = -13890789
$ws.Range("H76").Value = 3741.1667
$ws.Range("I76").Value = 3708.8333
$ws.Range("J76").Value = 3757.3333
$ws.Range("K76").Value = 3708.8333
$ws.Range("L76").Value = 3757.3333
$ws.Range("M76").Value = -3393.8333
$ws.Range("N76").Value = -4387.3333
$ws.Range("H79").Value = 3741.1667
$ws.Range("I79").Value = 3708.8333
$ws.Range("J79").Value = 3757.3333
$ws.Range("K79").Value = 3708.8333
$ws.Range("L79").Value = 3757.3333
$ws.Range("M79").Value = -2616.8333
$ws.Range("N79").Value = -5941.3333
$ws.Range("H107").Value = 7564.8335
$ws.Range("I107").Value = 7077.8
$ws.Range("K107").Value = 7077.8
$ws.Range("M107").Value = -5157.8
$ws.Range("H116").Value = 2364.2856
$ws.Range("I116").Value = 2445.6428
$ws.Range("K116").Value = 2445.6428
$ws.Range("M116").Value = 996.3571999999999
$ws.Range("H132").Value = 15884589
$ws.Range("I132").Value = 18528354
$ws.Range("J132").Value = 22004
$ws.Range("K132").Value = 55585062
$ws.Range("L132").Value = 66012
$ws.Range("M132").Value = -55582532
$ws.Range("N132").Value = -71072
$ws.Range("H137").Value = 2358.9092
$ws.Range("I137").Value = 2023.3158
$ws.Range("J137").Value = 2814.3572
$ws.Range("K137").Value = 6069.9474
$ws.Range("L137").Value = 8443.071599999999
$ws.Range("M137").Value = -3519.9474
$ws.Range("N137").Value = -13543.0716
$ws.Range("H138").Value = 2679.0254
$ws.Range("I138").Value = 2308.4546
$ws.Range("J138").Value = 2738.9707
$ws.Range("K138").Value = 6925.3638
$ws.Range("L138").Value = 8216.9121
$ws.Range("M138").Value = -1785.3638
$ws.Range("N138").Value = -18496.9121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10847.189
$ws.Range("I32").Value = 8080.7207
$ws.Range("K32").Value = 8080.7207
$ws.Range("M32").Value = -7793.7207
$ws.Range("H63").Value = 2030.25
$ws.Range("I63").Value = 1900.16
$ws.Range("J63").Value = 2494.8572
$ws.Range("K63").Value = 1900.16
$ws.Range("L63").Value = 2494.8572
$ws.Range("M63").Value = -1214.16
$ws.Range("N63").Value = -3866.8572
$ws.Range("H66").Value = 2030.25
$ws.Range("I66").Value = 1900.16
$ws.Range("J66").Value = 2494.8572
$ws.Range("K66").Value = 9500.800000000001
$ws.Range("L66").Value = 12474.286
$ws.Range("M66").Value = -6068.800000000001
$ws.Range("N66").Value = -19338.286
$ws.Range("H74").Value = 3167.2144
$ws.Range("I74").Value = 1618.5714
$ws.Range("J74").Value = 4715.857
$ws.Range("K74").Value = 1618.5714
$ws.Range("L74").Value = 4715.857
$ws.Range("M74").Value = -744.5714
$ws.Range("N74").Value = -6463.857
$ws.Range("H77").Value = 3167.2144
$ws.Range("I77").Value = 1618.5714
$ws.Range("J77").Value = 4715.857
$ws.Range("K77").Value = 8092.857
$ws.Range("L77").Value = 23579.285
$ws.Range("M77").Value = -3724.857
$ws.Range("N77").Value = -32315.285
$ws.Range("H102").Value = 27782330
$ws.Range("I102").Value = 41668496
$ws.Range("K102").Value = 41668496
$ws.Range("M102").Value = -41666874
$ws.Range("H110").Value = 947.25
$ws.Range("I110").Value = 929.6667
$ws.Range("K110").Value = 929.6667
$ws.Range("M110").Value = 1115.3333
$ws.Range("H132").Value = 2426.426
$ws.Range("I132").Value = 1915.0714
$ws.Range("J132").Value = 4216.1665
$ws.Range("K132").Value = 5745.2142
$ws.Range("L132").Value = 12648.4995
$ws.Range("M132").Value = -3215.2142
$ws.Range("N132").Value = -17708.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 293.44446
$ws.Range("I22").Value = 255.125
$ws.Range("K22").Value = 255.125
$ws.Range("M22").Value = -82.125
$ws.Range("H99").Value = 55556692
$ws.Range("I99").Value = 83334380
$ws.Range("J99").Value = 1303.5
$ws.Range("K99").Value = 83334380
$ws.Range("L99").Value = 1303.5
$ws.Range("M99").Value = -83332882
$ws.Range("N99").Value = -4299.5
$ws.Range("H134").Value = 6201.75
$ws.Range("I134").Value = 970.4
$ws.Range("K134").Value = 2911.2
$ws.Range("M134").Value = -376.1999999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1396.7858
$ws.Range("I31").Value = 1352.4755
$ws.Range("K31").Value = 1352.4755
$ws.Range("M31").Value = -1057.4755
$ws.Range("H34").Value = 1396.7858
$ws.Range("I34").Value = 1352.4755
$ws.Range("K34").Value = 1352.4755
$ws.Range("M34").Value = -1150.4755
$ws.Range("H132").Value = 1820.919
$ws.Range("I132").Value = 1471.6086
$ws.Range("J132").Value = 2394.7856
$ws.Range("K132").Value = 4414.825800000001
$ws.Range("L132").Value = 7184.3568
$ws.Range("M132").Value = -1884.825800000001
$ws.Range("N132").Value = -12244.3568
$ws.Range("H134").Value = 9616612
$ws.Range("I134").Value = 1235.5526
$ws.Range("J134").Value = 35715490
$ws.Range("K134").Value = 3706.6578
$ws.Range("L134").Value = 107146470
$ws.Range("M134").Value = -1171.6578
$ws.Range("N134").Value = -107151540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12591.2
$ws.Range("I3").Value = 7429
$ws.Range("J3").Value = 18900.555
$ws.Range("K3").Value = 22287
$ws.Range("L3").Value = 56701.665
$ws.Range("M3").Value = -22175
$ws.Range("N3").Value = -56925.665
$ws.Range("H93").Value = 6424.5
$ws.Range("J93").Value = 6424.5
$ws.Range("L93").Value = 19273.5
$ws.Range("N93").Value = -23017.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20518
$ws.Range("H80").Value = 3611.5625
$ws.Range("I80").Value = 2161.6667
$ws.Range("J80").Value = 3946.1538
$ws.Range("K80").Value = 2161.6667
$ws.Range("L80").Value = 3946.1538
$ws.Range("M80").Value = -1163.6667
$ws.Range("N80").Value = -5942.1538
$ws.Range("H83").Value = 3611.5625
$ws.Range("I83").Value = 2161.6667
$ws.Range("J83").Value = 3946.1538
$ws.Range("K83").Value = 10808.3335
$ws.Range("L83").Value = 19730.769
$ws.Range("M83").Value = -5816.333500000001
$ws.Range("N83").Value = -29714.769
$ws.Range("H102").Value = 2251.1516
$ws.Range("I102").Value = 1419.3914
$ws.Range("K102").Value = 1419.3914
$ws.Range("M102").Value = 202.6086
$ws.Range("H107").Value = 88.85714
$ws.Range("I107").Value = 93
$ws.Range("J107").Value = 83.333336
$ws.Range("K107").Value = 93
$ws.Range("L107").Value = 83.333336
$ws.Range("M107").Value = 1827
$ws.Range("N107").Value = -3923.333336
$ws.Range("H113").Value = 1737.1428
$ws.Range("I113").Value = 1653.3334
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1653.3334
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 516.6666
$ws.Range("N113").Value = -6140
$ws.Range("H132").Value = 6512.393
$ws.Range("I132").Value = 8113
$ws.Range("J132").Value = 4038.7273
$ws.Range("K132").Value = 24339
$ws.Range("L132").Value = 12116.1819
$ws.Range("M132").Value = -21809
$ws.Range("N132").Value = -17176.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1880.8182
$ws.Range("I7").Value = 1523
$ws.Range("J7").Value = 2835
$ws.Range("K7").Value = 1523
$ws.Range("L7").Value = 2835
$ws.Range("M7").Value = -1411
$ws.Range("N7").Value = -3059
$ws.Range("H16").Value = 1414.3846
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 5744.5557
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5744.5557
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5744.5557
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -6120.5557
$ws.Range("H61").Value = 1163.3572
$ws.Range("I61").Value = 1163.3572
$ws.Range("K61").Value = 1163.3572
$ws.Range("M61").Value = -961.3571999999999
$ws.Range("H68").Value = 1302
$ws.Range("I68").Value = 1302
$ws.Range("K68").Value = 1302
$ws.Range("M68").Value = -553
$ws.Range("H71").Value = 1302
$ws.Range("I71").Value = 1302
$ws.Range("K71").Value = 6510
$ws.Range("M71").Value = -2766
$ws.Range("H82").Value = 1037.2222
$ws.Range("I82").Value = 1022.5
$ws.Range("J82").Value = 1041.4286
$ws.Range("K82").Value = 1022.5
$ws.Range("L82").Value = 1041.4286
$ws.Range("M82").Value = -661.5
$ws.Range("N82").Value = -1763.4286
$ws.Range("H85").Value = 1037.2222
$ws.Range("I85").Value = 1022.5
$ws.Range("J85").Value = 1041.4286
$ws.Range("K85").Value = 1022.5
$ws.Range("L85").Value = 1041.4286
$ws.Range("M85").Value = 225.5
$ws.Range("N85").Value = -3537.4286
$ws.Range("H107").Value = 19565.5
$ws.Range("I107").Value = 19565.5
$ws.Range("K107").Value = 19565.5
$ws.Range("M107").Value = -17645.5
$ws.Range("H113").Value = 1163.3572
$ws.Range("I113").Value = 1163.3572
$ws.Range("K113").Value = 1163.3572
$ws.Range("M113").Value = 1006.6428
$ws.Range("H126").Value = 1880.8182
$ws.Range("I126").Value = 1523
$ws.Range("J126").Value = 2835
$ws.Range("K126").Value = 4569
$ws.Range("L126").Value = 8505
$ws.Range("M126").Value = -2099
$ws.Range("N126").Value = -13445
$ws.Range("H132").Value = 2393.1516
$ws.Range("I132").Value = 1958.7368
$ws.Range("J132").Value = 2982.7144
$ws.Range("K132").Value = 5876.2104
$ws.Range("L132").Value = 8948.143199999999
$ws.Range("M132").Value = -3346.2104
$ws.Range("N132").Value = -14008.1432
$ws.Range("H141").Value = 59621.668
$ws.Range("J141").Value = 59546
$ws.Range("L141").Value = 59546
$ws.Range("N141").Value = -69906

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3731.4146
$ws.Range("I132").Value = 4037.0312
$ws.Range("J132").Value = 2644.7778
$ws.Range("K132").Value = 12111.0936
$ws.Range("L132").Value = 7934.3334
$ws.Range("M132").Value = -9581.0936
$ws.Range("N132").Value = -12994.3334
$ws.Range("H136").Value = 1625.4318
$ws.Range("I136").Value = 735.2105
$ws.Range("K136").Value = 2205.6315
$ws.Range("M136").Value = 344.3685

Write-Host "Applied market price updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
